# Auto-generated edit script applying the diff to Malboro_Profits workbook
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 1185.091  # H9: '1089.4166' -> '1185.091'
$ws.Cells.Item(9, 9).Value = 1292.6666  # I9: '1167.1' -> '1292.6666'
$ws.Cells.Item(9, 11).Value = 1292.6666  # K9: '1167.1' -> '1292.6666'
$ws.Cells.Item(9, 13).Value = -1123.6666  # M9: '-998.0999999999999' -> '-1123.6666'
$ws.Cells.Item(19, 8).Value = 70828.164  # H19: '84793.8' -> '70828.164'
$ws.Cells.Item(19, 10).Value = 84813.8  # J19: '105767.25' -> '84813.8'
$ws.Cells.Item(19, 12).Value = 84813.8  # L19: '105767.25' -> '84813.8'
$ws.Cells.Item(19, 14).Value = -85163.8  # N19: '-106117.25' -> '-85163.8'
$ws.Cells.Item(32, 8).Value = 1856  # H32: '1532.8572' -> '1856'
$ws.Cells.Item(32, 9).Value = 2280  # I32: '1490' -> '2280'
$ws.Cells.Item(32, 10).Value = 1750  # J32: '1550' -> '1750'
$ws.Cells.Item(32, 11).Value = 2280  # K32: '1490' -> '2280'
$ws.Cells.Item(32, 12).Value = 1750  # L32: '1550' -> '1750'
$ws.Cells.Item(32, 13).Value = -1954  # M32: '-1164' -> '-1954'
$ws.Cells.Item(32, 14).Value = -2402  # N32: '-2202' -> '-2402'
$ws.Cells.Item(33, 8).Value = 13160347  # H33: '12502627' -> '13160347'
$ws.Cells.Item(33, 9).Value = 25000860  # I33: '25000834' -> '25000860'
$ws.Cells.Item(33, 10).Value = 4222.3335  # J33: '4420.1' -> '4222.3335'
$ws.Cells.Item(33, 11).Value = 25000860  # K33: '25000834' -> '25000860'
$ws.Cells.Item(33, 12).Value = 4222.3335  # L33: '4420.1' -> '4222.3335'
$ws.Cells.Item(33, 13).Value = -25000631  # M33: '-25000605' -> '-25000631'
$ws.Cells.Item(33, 14).Value = -4680.3335  # N33: '-4878.1' -> '-4680.3335'
$ws.Cells.Item(64, 8).Value = 7102.0557  # H64: '6991.421' -> '7102.0557'
$ws.Cells.Item(64, 9).Value = 3424.5  # I64: '3949.6667' -> '3424.5'
$ws.Cells.Item(64, 11).Value = 3424.5  # K64: '3949.6667' -> '3424.5'
$ws.Cells.Item(64, 13).Value = -3176.5  # M64: '-3701.6667' -> '-3176.5'
$ws.Cells.Item(67, 8).Value = 7102.0557  # H67: '6991.421' -> '7102.0557'
$ws.Cells.Item(67, 9).Value = 3424.5  # I67: '3949.6667' -> '3424.5'
$ws.Cells.Item(67, 11).Value = 3424.5  # K67: '3949.6667' -> '3424.5'
$ws.Cells.Item(67, 13).Value = -2566.5  # M67: '-3091.6667' -> '-2566.5'
$ws.Cells.Item(106, 8).Value = 5848.3335  # H106: '5848.1665' -> '5848.3335'
$ws.Cells.Item(106, 9).Value = 5848.3335  # I106: '6718' -> '5848.3335'
$ws.Cells.Item(106, 10).Value = 0  # J106: '1499' -> '0'
$ws.Cells.Item(106, 11).Value = 5848.3335  # K106: '6718' -> '5848.3335'
$ws.Cells.Item(106, 12).Value = 0  # L106: '1499' -> '0'
$ws.Cells.Item(106, 13).Value = -5217.3335  # M106: '-6087' -> '-5217.3335'
$ws.Cells.Item(106, 14).ClearContents()  # N106: was '-2761'
$ws.Cells.Item(127, 8).Value = 839443.3  # H127: '719908.5600000001' -> '839443.3'
$ws.Cells.Item(127, 9).Value = 1670220  # I127: '1252831.5' -> '1670220'
$ws.Cells.Item(127, 10).Value = 8666.666999999999  # J127: '9344.666999999999' -> '8666.666999999999'
$ws.Cells.Item(127, 11).Value = 5010660  # K127: '3758494.5' -> '5010660'
$ws.Cells.Item(127, 12).Value = 26000.001  # L127: '28034.001' -> '26000.001'
$ws.Cells.Item(127, 13).Value = -5005700  # M127: '-3753534.5' -> '-5005700'
$ws.Cells.Item(127, 14).Value = -35920.001  # N127: '-37954.001' -> '-35920.001'
$ws.Cells.Item(135, 8).Value = 1576.1111  # H135: '1610.3529' -> '1576.1111'
$ws.Cells.Item(135, 9).Value = 1691.5333  # I135: '1691.8' -> '1691.5333'
$ws.Cells.Item(135, 10).Value = 999  # J135: '999.5' -> '999'
$ws.Cells.Item(135, 11).Value = 15223.7997  # K135: '15226.2' -> '15223.7997'
$ws.Cells.Item(135, 12).Value = 8991  # L135: '8995.5' -> '8991'
$ws.Cells.Item(135, 13).Value = -12688.7997  # M135: '-12691.2' -> '-12688.7997'
$ws.Cells.Item(135, 14).Value = -14061  # N135: '-14065.5' -> '-14061'

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 2007  # H26: '737.1429000000001' -> '2007'
$ws.Cells.Item(26, 9).Value = 2007  # I26: '737.1429000000001' -> '2007'
$ws.Cells.Item(26, 11).Value = 2007  # K26: '737.1429000000001' -> '2007'
$ws.Cells.Item(26, 13).Value = -1677  # M26: '-407.1429000000001' -> '-1677'
$ws.Cells.Item(74, 8).Value = 14296.718  # H74: '13957.1' -> '14296.718'
$ws.Cells.Item(74, 10).Value = 30485.47  # J74: '28831.389' -> '30485.47'
$ws.Cells.Item(74, 12).Value = 30485.47  # L74: '28831.389' -> '30485.47'
$ws.Cells.Item(74, 14).Value = -32233.47  # N74: '-30579.389' -> '-32233.47'
$ws.Cells.Item(77, 8).Value = 14296.718  # H77: '13957.1' -> '14296.718'
$ws.Cells.Item(77, 10).Value = 30485.47  # J77: '28831.389' -> '30485.47'
$ws.Cells.Item(77, 12).Value = 152427.35  # L77: '144156.945' -> '152427.35'
$ws.Cells.Item(77, 14).Value = -161163.35  # N77: '-152892.945' -> '-161163.35'
$ws.Cells.Item(110, 8).Value = 25899  # H110: '49998' -> '25899'
$ws.Cells.Item(110, 9).Value = 25899  # I110: '49998' -> '25899'
$ws.Cells.Item(110, 11).Value = 25899  # K110: '49998' -> '25899'
$ws.Cells.Item(110, 13).Value = -23854  # M110: '-47953' -> '-23854'
$ws.Cells.Item(122, 8).Value = 2626  # H122: '2862.7' -> '2626'
$ws.Cells.Item(122, 9).Value = 2546.5454  # I122: '2791.889' -> '2546.5454'
$ws.Cells.Item(122, 11).Value = 7639.6362  # K122: '8375.667000000001' -> '7639.6362'
$ws.Cells.Item(122, 13).Value = -5189.6362  # M122: '-5925.667000000001' -> '-5189.6362'

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 9523.625  # H99: '3219.2068' -> '9523.625'
$ws.Cells.Item(99, 9).Value = 10227  # I99: '3023.25' -> '10227'
$ws.Cells.Item(99, 10).Value = 4600  # J99: '4159.8' -> '4600'
$ws.Cells.Item(99, 11).Value = 10227  # K99: '3023.25' -> '10227'
$ws.Cells.Item(99, 12).Value = 4600  # L99: '4159.8' -> '4600'
$ws.Cells.Item(99, 13).Value = -8729  # M99: '-1525.25' -> '-8729'
$ws.Cells.Item(99, 14).Value = -7596  # N99: '-7155.8' -> '-7596'
$ws.Cells.Item(107, 8).Value = 1969.0605  # H107: '2072.9333' -> '1969.0605'
$ws.Cells.Item(107, 9).Value = 2133.8696  # I107: '2247.0476' -> '2133.8696'
$ws.Cells.Item(107, 10).Value = 1590  # J107: '1666.6666' -> '1590'
$ws.Cells.Item(107, 11).Value = 2133.8696  # K107: '2247.0476' -> '2133.8696'
$ws.Cells.Item(107, 12).Value = 1590  # L107: '1666.6666' -> '1590'
$ws.Cells.Item(107, 13).Value = -213.8696  # M107: '-327.0475999999999' -> '-213.8696'
$ws.Cells.Item(107, 14).Value = -5430  # N107: '-5506.6666' -> '-5430'
$ws.Cells.Item(134, 8).Value = 41501.03  # H134: '42818.805' -> '41501.03'
$ws.Cells.Item(134, 9).Value = 49138.953  # I134: '51563.4' -> '49138.953'
$ws.Cells.Item(134, 11).Value = 147416.859  # K134: '154690.2' -> '147416.859'
$ws.Cells.Item(134, 13).Value = -144881.859  # M134: '-152155.2' -> '-144881.859'

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 9682.6  # H86: '10375.111' -> '9682.6'
$ws.Cells.Item(86, 10).Value = 5195.5557  # J86: '5694.2856' -> '5195.5557'
$ws.Cells.Item(86, 12).Value = 5195.5557  # L86: '5694.2856' -> '5195.5557'
$ws.Cells.Item(86, 14).Value = -7441.5557  # N86: '-7940.2856' -> '-7441.5557'
$ws.Cells.Item(89, 8).Value = 9682.6  # H89: '10375.111' -> '9682.6'
$ws.Cells.Item(89, 10).Value = 5195.5557  # J89: '5694.2856' -> '5195.5557'
$ws.Cells.Item(89, 12).Value = 25977.7785  # L89: '28471.428' -> '25977.7785'
$ws.Cells.Item(89, 14).Value = -37209.7785  # N89: '-39703.428' -> '-37209.7785'
$ws.Cells.Item(94, 8).Value = 883.6667  # H94: '711.7778' -> '883.6667'
$ws.Cells.Item(94, 9).Value = 0  # I94: '330.66666' -> '0'
$ws.Cells.Item(94, 10).Value = 883.6667  # J94: '902.3333' -> '883.6667'
$ws.Cells.Item(94, 11).Value = 0  # K94: '330.66666' -> '0'
$ws.Cells.Item(94, 12).Value = 883.6667  # L94: '902.3333' -> '883.6667'
$ws.Cells.Item(94, 13).ClearContents()  # M94: was '120.33334'
$ws.Cells.Item(94, 14).Value = -1785.6667  # N94: '-1804.3333' -> '-1785.6667'
$ws.Cells.Item(134, 8).Value = 21281428  # H134: '22732340' -> '21281428'
$ws.Cells.Item(134, 9).Value = 2219.9688  # I134: '2334.2666' -> '2219.9688'
$ws.Cells.Item(134, 10).Value = 66677068  # J134: '71439496' -> '66677068'
$ws.Cells.Item(134, 11).Value = 6659.9064  # K134: '7002.7998' -> '6659.9064'
$ws.Cells.Item(134, 12).Value = 200031204  # L134: '214318488' -> '200031204'
$ws.Cells.Item(134, 13).Value = -4124.9064  # M134: '-4467.7998' -> '-4124.9064'
$ws.Cells.Item(134, 14).Value = -200036274  # N134: '-214323558' -> '-200036274'
$ws.Cells.Item(135, 8).Value = 92000  # H135: '0' -> '92000'
$ws.Cells.Item(135, 10).Value = 92000  # J135: '0' -> '92000'
$ws.Cells.Item(135, 12).Value = 92000  # L135: '0' -> '92000'
$ws.Cells.Item(135, 14).Value = -102140  # N135: None -> '-102140'

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 872.1111  # H68: '979.625' -> '872.1111'
$ws.Cells.Item(68, 10).Value = 1062.3334  # J68: '1272.4' -> '1062.3334'
$ws.Cells.Item(68, 12).Value = 3187.0002  # L68: '3817.2' -> '3187.0002'
$ws.Cells.Item(68, 14).Value = -4809.0002  # N68: '-5439.200000000001' -> '-4809.0002'
$ws.Cells.Item(71, 8).Value = 872.1111  # H71: '979.625' -> '872.1111'
$ws.Cells.Item(71, 10).Value = 1062.3334  # J71: '1272.4' -> '1062.3334'
$ws.Cells.Item(71, 12).Value = 9561.000599999999  # L71: '11451.6' -> '9561.000599999999'
$ws.Cells.Item(71, 14).Value = -17673.0006  # N71: '-19563.6' -> '-17673.0006'
$ws.Cells.Item(86, 8).Value = 756.125  # H86: '631.4545000000001' -> '756.125'
$ws.Cells.Item(86, 9).Value = 723  # I86: '589.1053000000001' -> '723'
$ws.Cells.Item(86, 11).Value = 2169  # K86: '1767.3159' -> '2169'
$ws.Cells.Item(86, 13).Value = -983  # M86: '-581.3159000000001' -> '-983'
$ws.Cells.Item(89, 8).Value = 756.125  # H89: '631.4545000000001' -> '756.125'
$ws.Cells.Item(89, 9).Value = 723  # I89: '589.1053000000001' -> '723'
$ws.Cells.Item(89, 11).Value = 6507  # K89: '5301.947700000001' -> '6507'
$ws.Cells.Item(89, 13).Value = -579  # M89: '626.0522999999994' -> '-579'
$ws.Cells.Item(92, 8).Value = 321.18182  # H92: '308.75' -> '321.18182'
$ws.Cells.Item(92, 10).Value = 288.4  # J92: '269' -> '288.4'
$ws.Cells.Item(92, 12).Value = 865.1999999999999  # L92: '807' -> '865.1999999999999'
$ws.Cells.Item(92, 14).Value = -3361.2  # N92: '-3303' -> '-3361.2'
$ws.Cells.Item(139, 8).Value = 11452.9375  # H139: '11613.533' -> '11452.9375'
$ws.Cells.Item(139, 9).Value = 12374.786  # I139: '12631' -> '12374.786'
$ws.Cells.Item(139, 11).Value = 37124.358  # K139: '37893' -> '37124.358'
$ws.Cells.Item(139, 13).Value = -31984.358  # M139: '-32753' -> '-31984.358'

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 801.8333  # H107: '764.4737' -> '801.8333'
$ws.Cells.Item(107, 9).Value = 860.5625  # I107: '815.35297' -> '860.5625'
$ws.Cells.Item(107, 11).Value = 860.5625  # K107: '815.35297' -> '860.5625'
$ws.Cells.Item(107, 13).Value = 1059.4375  # M107: '1104.64703' -> '1059.4375'
$ws.Cells.Item(122, 8).Value = 7151.8  # H122: '7126.4165' -> '7151.8'
$ws.Cells.Item(122, 9).Value = 8216  # I122: '7945.6665' -> '8216'
$ws.Cells.Item(122, 11).Value = 24648  # K122: '23836.9995' -> '24648'
$ws.Cells.Item(122, 13).Value = -22198  # M122: '-21386.9995' -> '-22198'
$ws.Cells.Item(123, 8).Value = 58724  # H123: '58774.25' -> '58724'
$ws.Cells.Item(123, 10).Value = 58724  # J123: '58774.25' -> '58724'
$ws.Cells.Item(123, 12).Value = 58724  # L123: '58774.25' -> '58724'
$ws.Cells.Item(123, 14).Value = -63624  # N123: '-63674.25' -> '-63624'

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(19, 8).Value = 6000  # H19: '0' -> '6000'
$ws.Cells.Item(19, 9).Value = 2000  # I19: '0' -> '2000'
$ws.Cells.Item(19, 10).Value = 10000  # J19: '0' -> '10000'
$ws.Cells.Item(19, 11).Value = 2000  # K19: '0' -> '2000'
$ws.Cells.Item(19, 12).Value = 10000  # L19: '0' -> '10000'
$ws.Cells.Item(19, 13).Value = -1830  # M19: None -> '-1830'
$ws.Cells.Item(19, 14).Value = -10340  # N19: None -> '-10340'
$ws.Cells.Item(24, 8).Value = 66777.78  # H24: '75000' -> '66777.78'
$ws.Cells.Item(24, 10).Value = 80142.86  # J24: '93333.336' -> '80142.86'
$ws.Cells.Item(24, 12).Value = 80142.86  # L24: '93333.336' -> '80142.86'
$ws.Cells.Item(24, 14).Value = -80828.86  # N24: '-94019.336' -> '-80828.86'
$ws.Cells.Item(25, 8).Value = 212249.75  # H25: '172799.8' -> '212249.75'
$ws.Cells.Item(25, 10).Value = 278333.34  # J25: '212500' -> '278333.34'
$ws.Cells.Item(25, 12).Value = 278333.34  # L25: '212500' -> '278333.34'
$ws.Cells.Item(25, 14).Value = -278793.34  # N25: '-212960' -> '-278793.34'
$ws.Cells.Item(50, 8).Value = 58841.5  # H50: '53893.332' -> '58841.5'
$ws.Cells.Item(50, 10).Value = 58841.5  # J50: '53893.332' -> '58841.5'
$ws.Cells.Item(50, 12).Value = 58841.5  # L50: '53893.332' -> '58841.5'
$ws.Cells.Item(50, 14).Value = -60115.5  # N50: '-55167.332' -> '-60115.5'
$ws.Cells.Item(54, 8).Value = 84210  # H54: '74980' -> '84210'
$ws.Cells.Item(54, 10).Value = 84210  # J54: '74980' -> '84210'
$ws.Cells.Item(54, 12).Value = 84210  # L54: '74980' -> '84210'
$ws.Cells.Item(54, 14).Value = -85498  # N54: '-76268' -> '-85498'
$ws.Cells.Item(55, 8).Value = 2171.3333  # H55: '2278.9412' -> '2171.3333'
$ws.Cells.Item(55, 9).Value = 1903.8889  # I55: '2342.1428' -> '1903.8889'
$ws.Cells.Item(55, 10).Value = 2438.7778  # J55: '2234.7' -> '2438.7778'
$ws.Cells.Item(55, 11).Value = 1903.8889  # K55: '2342.1428' -> '1903.8889'
$ws.Cells.Item(55, 12).Value = 2438.7778  # L55: '2234.7' -> '2438.7778'
$ws.Cells.Item(55, 13).Value = -1730.8889  # M55: '-2169.1428' -> '-1730.8889'
$ws.Cells.Item(55, 14).Value = -2784.7778  # N55: '-2580.7' -> '-2784.7778'
$ws.Cells.Item(93, 8).Value = 4304.615  # H93: '4459.04' -> '4304.615'
$ws.Cells.Item(93, 9).Value = 3962.1  # I93: '4147.263' -> '3962.1'
$ws.Cells.Item(93, 11).Value = 3962.1  # K93: '4147.263' -> '3962.1'
$ws.Cells.Item(93, 13).Value = -2714.1  # M93: '-2899.263' -> '-2714.1'
$ws.Cells.Item(136, 8).Value = 112925.62  # H136: '124517.63' -> '112925.62'
$ws.Cells.Item(136, 9).Value = 16962.467  # I136: '18059.5' -> '16962.467'
$ws.Cells.Item(136, 10).Value = 352833.5  # J136: '422600.4' -> '352833.5'
$ws.Cells.Item(136, 11).Value = 50887.401  # K136: '54178.5' -> '50887.401'
$ws.Cells.Item(136, 12).Value = 1058500.5  # L136: '1267801.2' -> '1058500.5'
$ws.Cells.Item(136, 13).Value = -48337.401  # M136: '-51628.5' -> '-48337.401'
$ws.Cells.Item(136, 14).Value = -1063600.5  # N136: '-1272901.2' -> '-1063600.5'

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 19000  # H20: '14500' -> '19000'
$ws.Cells.Item(20, 9).Value = 0  # I20: '10000' -> '0'
$ws.Cells.Item(20, 11).Value = 0  # K20: '10000' -> '0'
$ws.Cells.Item(20, 13).ClearContents()  # M20: was '-9760'
$ws.Cells.Item(126, 8).Value = 5650  # H126: '20002738' -> '5650'
$ws.Cells.Item(126, 9).Value = 5650  # I126: '3422.5' -> '5650'
$ws.Cells.Item(126, 10).Value = 0  # J126: '100000000' -> '0'
$ws.Cells.Item(126, 11).Value = 16950  # K126: '10267.5' -> '16950'
$ws.Cells.Item(126, 12).Value = 0  # L126: '300000000' -> '0'
$ws.Cells.Item(126, 13).Value = -14480  # M126: '-7797.5' -> '-14480'
$ws.Cells.Item(126, 14).ClearContents()  # N126: was '-300004940'
$ws.Cells.Item(132, 8).Value = 5164.793  # H132: '5174.9653' -> '5164.793'
$ws.Cells.Item(132, 9).Value = 1843.2  # I132: '1855' -> '1843.2'
$ws.Cells.Item(132, 11).Value = 5529.6  # K132: '5565' -> '5529.6'
$ws.Cells.Item(132, 13).Value = -2999.6  # M132: '-3035' -> '-2999.6'
